# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# The G column (header "K") previously held a "Strike#" count; it is being
# regenerated to hold actual strikeout (K) values. Only column G (rows 2-44)
# changes; every other column is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 10
    3  = 7
    4  = 6
    5  = 10
    6  = 5
    7  = 5
    8  = 3
    9  = 2
    10 = 8
    11 = 7
    12 = 8
    13 = 2
    14 = 10
    15 = 5
    16 = 1
    17 = 13
    18 = 8
    19 = 6
    20 = 7
    21 = 0
    22 = 4
    23 = 0
    24 = 0
    25 = 0
    26 = 1
    27 = 3
    28 = 1
    29 = 3
    30 = 3
    31 = 3
    32 = 6
    33 = 2
    34 = 6
    35 = 5
    36 = 7
    37 = 4
    38 = 6
    39 = 6
    40 = 5
    41 = 5
    42 = 5
    43 = 2
    44 = 4
}

foreach ($row in $newK.Keys) {
    $ws.Cells.Item($row, 7).Value = $newK[$row]
}
